$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row formatting (A1:B1) onto the destination header cells
# (E4:F4) before the source is cleared, so the existing named/shared style
# (bold, border, center/top aligned) carries over instead of creating new
# style entries.
$ws.Range("A1:B1").Copy()
$ws.Range("E4:F4").PasteSpecial(-4122)  # xlPasteFormats

# Clear the old data range A1:B5
$ws.Range("A1:B5").Clear()

# Write the headers and data to the new location E4:F8
$ws.Range("E4").Value = "Name"
$ws.Range("F4").Value = "Height"

$ws.Range("E5").Value = "Adiya"
$ws.Range("F5").Value = 179

$ws.Range("E6").Value = "Samen"
$ws.Range("F6").Value = 181

$ws.Range("E7").Value = "Darek"
$ws.Range("F7").Value = 170

$ws.Range("E8").Value = "Jan"
$ws.Range("F8").Value = 167
